# Update cryptocurrency price ("Price", column D) and 1-hour change
# ("Volume(1h)", column E) figures to the latest scraped values.
#
# For column D cells whose new text would otherwise be auto-parsed by
# Excel as a number (dropping e.g. a trailing ".00"), the cell is
# temporarily switched to text format ("@"), the literal string is
# written, and the style is then reset back to "Normal" so the cell
# ends up with the same (default) style it started with while still
# holding the exact text we want.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.383.94"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "3.782.79"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -10.10%  "
$ws.Range("D7").Value = "3.779.04"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.716"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.32%  "
$ws.Range("E11").Value = "  -9.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -11.50%  "
$ws.Range("D14").Value = "4.397.85"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +17.90%  "
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "3.784.01"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.02%  "
$ws.Range("D20").Value = "66.528.85"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "401.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.59%  "
$ws.Range("E28").Value = "  -4.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "700.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -7.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("D39").Value = "0.0₃0761"
$ws.Range("E39").Value = "  +10.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0447"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.05%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.133"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.62%  "
